$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "[Diana%Aguiar de Sousa%NULL%1,         H Bart%van der Worp%NULL%2,         H Bart%van der Worp%NULL%0,         Valeria%Caso%NULL%1,         Charlotte%Cordonnier%NULL%1,         Daniel%Strbian%NULL%2,         Daniel%Strbian%NULL%0,         George%Ntaios%NULL%1,         Peter D%Schellinger%NULL%1,         Else Charlotte%Sandset%NULL%1,         NULL%NULL%NULL%0,         NULL%NULL%NULL%0]"
$ws.Range("E3").Value = "[Yogini V.%Chudasama%NULL%1,         Clare L.%Gillies%NULL%1,         Francesco%Zaccardi%NULL%0,         Briana%Coles%NULL%1,         Melanie J.%Davies%NULL%1,         Samuel%Seidu%NULL%1,         Kamlesh%Khunti%NULL%0]"
$ws.Range("E4").Value = "[Rocco%Capuano%NULL%1,         Manuela%Altieri%NULL%1,         Alvino%Bisecco%NULL%1,         Alessandro%d’Ambrosio%NULL%1,         Renato%Docimo%NULL%1,         Daniela%Buonanno%NULL%1,         Federica%Matrone%NULL%1,         Federica%Giuliano%NULL%1,         Gioacchino%Tedeschi%NULL%1,         Gabriella%Santangelo%NULL%1,         Antonio%Gallo%antonio.gallo@unicampania.it%1]"
$ws.Range("E5").Value = "[Jing%Zhao%NULL%0,         Hang%Li%NULL%1,         David%Kung%NULL%1,         Marc%Fisher%NULL%1,         Ying%Shen%NULL%1,         Renyu%Liu%NULL%1]"
$ws.Range("E6").Value = "[Farrah J.%Mateen%fmateen@mgh.harvard.edu%1,         Shawheen%Rezaei%NULL%2,         Shawheen%Rezaei%NULL%0,         Nicholas%Alakel%NULL%1,         Brittany%Gazdag%NULL%1,         Aditya Ravi%Kumar%NULL%1,         Andre%Vogel%NULL%1]"
$ws.Range("E7").Value = "[Ljiljana%Radulovic%NULL%1,         Jevto%Erakovic%NULL%1,         Milovan%Roganovic%NULL%1]"
$ws.Range("E8").Value = "[Aleksandar%Stojanov%NULL%1,         Marina%Malobabic%NULL%1,         Vuk%Milosevic%NULL%1,         Jelena%Stojanov%NULL%1,         Slobodan%Vojinovic%NULL%1,         Goran%Stanojevic%NULL%1,         Milos%Stevic%NULL%1]"
$ws.Range("E9").Value = "[Carla%Piano%NULL%1,         Enrico%Di Stasio%NULL%1,         Guido%Primiano%NULL%1,         Delfina%Janiri%NULL%1,         Marco%Luigetti%NULL%1,         Giovanni%Frisullo%NULL%1,         Catello%Vollono%NULL%1,         Matteo%Lucchini%NULL%1,         Valerio%Brunetti%NULL%1,         Mauro%Monforte%NULL%1,         Valeria%Guglielmi%NULL%1,         Giacomo%Della Marca%NULL%1,         Amelia%Evoli%NULL%1,         Camillo%Marra%NULL%1,         Massimiliano%Mirabella%NULL%1,         Davide%Quaranta%NULL%1,         Enzo%Ricci%NULL%1,         Serenella%Servidei%NULL%1,         Gabriella%Silvestri%NULL%1,         Simone%Bellavia%NULL%1,         Sara%Bortolani%NULL%1,         Francesco%Bove%NULL%1,         Riccardo%Di Iorio%NULL%1,         Andrea%Di Paolantonio%NULL%1,         Danilo%Genovese%NULL%1,         Tamara%Ialongo%NULL%1,         Maria Rita%Lo Monaco%NULL%1,         Jessica%Marotta%NULL%1,         Agata Katia%Patanella%NULL%1,         Alessia%Perna%NULL%1,         Martina%Petracca%NULL%1,         Giorgia%Presicce%NULL%1,         Vittorio%Riso%NULL%1,         Eleonora%Rollo%NULL%1,         Angela%Romano%NULL%1,         Marina%Romozzi%NULL%1,         Cristina%Sancricca%NULL%1,         Irene%Scala%NULL%1,         Gregorio%Spagni%NULL%1,         Marcella%Solito%NULL%1,         Luca%Tricoli%NULL%1,         Paola%Zinzi%NULL%1,         Paolo%Calabresi%NULL%1,         Anna Rita%Bentivoglio%NULL%1]"
$ws.Range("E10").Value = "[Accogli%M.%coreGivesNoEmail%1,       Barone%R.%coreGivesNoEmail%1,       Capozzo%R.%coreGivesNoEmail%1,       Logroscino%G.%coreGivesNoEmail%1,       Musio%M.%coreGivesNoEmail%1,       Zoccolella%S.%coreGivesNoEmail%1]"
$ws.Range("E11").Value = "[Vincenzo%Di Stefano%vincenzo19689@gmail.com%1,         Giuseppe%Battaglia%NULL%2,         Giuseppe%Battaglia%NULL%0,         Valerio%Giustino%NULL%1,         Andrea%Gagliardo%NULL%1,         Michele%D’Aleo%NULL%1,         Ottavio%Giannini%NULL%1,         Antonio%Palma%NULL%1,         Filippo%Brighina%NULL%1]"
$ws.Range("E12").Value = "[Guillermo F%López-Sánchez%gfls@um.es%1,         Rubén%López-Bueno%rlopezbu@unizar.es%1,         Alejandro%Gil-Salmerón%alejandro.gil.salmeron@uv.es%1,         Roksana%Zauder%NULL%1,         Maria%Skalska%NULL%1,         Joanna%Jastrzębska%NULL%1,         Zbigniew%Jastrzębski%NULL%1,         Felipe B%Schuch%NULL%1,         Igor%Grabovac%NULL%1,         Mark A%Tully%NULL%1,         Lee%Smith%NULL%0]"
